# Added pull vol script; added old update script; changed path to final output files
#
# - Make "trust" (sheet2) the active tab, and "bandit" (sheet1) no longer
#   selected.
# - Append 5 new values to the bottom of the "trust" sheet (A7:A11), which
#   extends the used range / dimension to A1:A11.
# - Update the active selection on "trust" to D22:D23.

$wb = $excel.ActiveWorkbook
$trust = $wb.Worksheets.Item("trust")

# Activating the sheet makes it the workbook's active tab (activeTab) and
# marks its sheetView as tabSelected, while clearing tabSelected on the
# previously-active sheet ("bandit").
$trust.Activate()

# Append the new rows below the existing data (A1:A6 already populated).
$trust.Range("A7").Value = 213163
$trust.Range("A8").Value = 221080
$trust.Range("A9").Value = 221099
$trust.Range("A10").Value = 202278
$trust.Range("A11").Value = 221183

# Move/extend the selection to match the final view state.
$trust.Range("D22:D23").Select()
